$wb = $excel.ActiveWorkbook

# The "True Positives" sheet (sheet1) gets a widened column G and
# cell G3 becomes the active/selected cell.
$ws = $wb.Worksheets.Item("True Positives")
$ws.Activate()

# Set width of column G (7th column) so the stored OOXML width is exactly 13
# (COM ColumnWidth and the stored <col width="..."/> differ by 5/MaxDigitWidth)
$ws.Columns.Item(7).ColumnWidth = 12.166666666666666

# Select G3 as the active cell
$ws.Range("G3").Select()
